# Auto-generated Excel COM-interop script
# Updates cached market-price columns (H-N) on multiple Leve sheets
# to refreshed values pulled by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 224779.22
$ws.Range("J17").Value = 251844.38
$ws.Range("L17").Value = 755533.14
$ws.Range("N17").Value = -755869.14
$ws.Range("H38").Value = 101.77778
$ws.Range("I38").Value = 101.77778
$ws.Range("K38").Value = 305.33334
$ws.Range("M38").Value = 66.66665999999998
$ws.Range("H43").Value = 1133.3334
$ws.Range("H58").Value = 3549.5833
$ws.Range("I58").Value = 307.66666
$ws.Range("J58").Value = 6791.5
$ws.Range("K58").Value = 922.9999799999999
$ws.Range("L58").Value = 20374.5
$ws.Range("M58").Value = -772.9999799999999
$ws.Range("N58").Value = -20674.5
$ws.Range("H76").Value = 3997787
$ws.Range("I76").Value = 4354.778
$ws.Range("K76").Value = 4354.778
$ws.Range("M76").Value = -4039.778
$ws.Range("H79").Value = 3997787
$ws.Range("I79").Value = 4354.778
$ws.Range("K79").Value = 4354.778
$ws.Range("M79").Value = -3262.778
$ws.Range("H80").Value = 1197134.8
$ws.Range("J80").Value = 492.6
$ws.Range("L80").Value = 1477.8
$ws.Range("N80").Value = -3473.8
$ws.Range("H83").Value = 1197134.8
$ws.Range("J83").Value = 492.6
$ws.Range("L83").Value = 4433.400000000001
$ws.Range("N83").Value = -14417.4
$ws.Range("H100").Value = 5426
$ws.Range("I100").Value = 5461.5
$ws.Range("K100").Value = 5461.5
$ws.Range("M100").Value = -4920.5
$ws.Range("H125").Value = 902.2857
$ws.Range("I125").Value = 913
$ws.Range("K125").Value = 8217
$ws.Range("M125").Value = -5757
$ws.Range("H137").Value = 3308.4846
$ws.Range("I137").Value = 2008.3889
$ws.Range("J137").Value = 3604.7087
$ws.Range("K137").Value = 6025.1667
$ws.Range("L137").Value = 10814.1261
$ws.Range("M137").Value = -3475.1667
$ws.Range("N137").Value = -15914.1261

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 165702.02
$ws.Range("I32").Value = 168413.72
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 168413.72
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -168126.72
$ws.Range("N32").Value = -3574
$ws.Range("H45").Value = 2347.7
$ws.Range("I45").Value = 1999.8334
$ws.Range("K45").Value = 1999.8334
$ws.Range("M45").Value = -1622.8334
$ws.Range("H46").Value = 9258.166999999999
$ws.Range("J46").Value = 9509.799999999999
$ws.Range("L46").Value = 9509.799999999999
$ws.Range("N46").Value = -10147.8
$ws.Range("H102").Value = 27394.607
$ws.Range("I102").Value = 6097.2383
$ws.Range("K102").Value = 6097.2383
$ws.Range("M102").Value = -4475.2383

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 19231080
$ws.Range("I94").Value = 20833502
$ws.Range("J94").Value = 2010
$ws.Range("K94").Value = 20833502
$ws.Range("L94").Value = 2010
$ws.Range("M94").Value = -20833051
$ws.Range("N94").Value = -2912
$ws.Range("H99").Value = 637.03845
$ws.Range("I99").Value = 616.8261
$ws.Range("J99").Value = 792
$ws.Range("K99").Value = 616.8261
$ws.Range("L99").Value = 792
$ws.Range("M99").Value = 881.1739
$ws.Range("N99").Value = -3788

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3966.87
$ws.Range("I31").Value = 1285.7142
$ws.Range("J31").Value = 4234.986
$ws.Range("K31").Value = 1285.7142
$ws.Range("L31").Value = 4234.986
$ws.Range("M31").Value = -990.7141999999999
$ws.Range("N31").Value = -4824.986
$ws.Range("H34").Value = 3966.87
$ws.Range("I34").Value = 1285.7142
$ws.Range("J34").Value = 4234.986
$ws.Range("K34").Value = 1285.7142
$ws.Range("L34").Value = 4234.986
$ws.Range("M34").Value = -1083.7142
$ws.Range("N34").Value = -4638.986
$ws.Range("H58").Value = 419307.25
$ws.Range("I58").Value = 1407.6
$ws.Range("K58").Value = 1407.6
$ws.Range("M58").Value = -1204.6
$ws.Range("H99").Value = 2674.5454
$ws.Range("I99").Value = 2008.1666
$ws.Range("K99").Value = 2008.1666
$ws.Range("M99").Value = -510.1666
$ws.Range("H126").Value = 2674.5454
$ws.Range("I126").Value = 2008.1666
$ws.Range("K126").Value = 6024.4998
$ws.Range("M126").Value = -3554.4998
$ws.Range("H130").Value = 40779.5
$ws.Range("J130").Value = 40779.5
$ws.Range("L130").Value = 40779.5
$ws.Range("N130").Value = -50819.5
$ws.Range("H136").Value = 419307.25
$ws.Range("I136").Value = 1407.6
$ws.Range("K136").Value = 4222.799999999999
$ws.Range("M136").Value = -1672.799999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 10666.417
$ws.Range("J98").Value = 1621.875
$ws.Range("L98").Value = 4865.625
$ws.Range("N98").Value = -7861.625
$ws.Range("H107").Value = 2012.4286
$ws.Range("I107").Value = 1865.1666
$ws.Range("J107").Value = 2122.875
$ws.Range("K107").Value = 5595.4998
$ws.Range("L107").Value = 6368.625
$ws.Range("M107").Value = -3675.4998
$ws.Range("N107").Value = -10208.625
$ws.Range("H124").Value = 1612.5
$ws.Range("I124").Value = 1612.5
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 4837.5
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = 72.5
$ws.Range("N124").ClearContents()
$ws.Range("H129").Value = 1624.1364
$ws.Range("I129").Value = 1121.4166
$ws.Range("J129").Value = 2227.4
$ws.Range("K129").Value = 3364.2498
$ws.Range("L129").Value = 6682.200000000001
$ws.Range("M129").Value = 1635.7502
$ws.Range("N129").Value = -16682.2
$ws.Range("H138").Value = 4074632.5
$ws.Range("I138").Value = 6001049
$ws.Range("K138").Value = 18003147
$ws.Range("M138").Value = -17998007

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 9848.267
$ws.Range("I126").Value = 16360.571
$ws.Range("K126").Value = 49081.713
$ws.Range("M126").Value = -46611.713

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3181
$ws.Range("I16").Value = 953
$ws.Range("J16").Value = 4666.3335
$ws.Range("K16").Value = 953
$ws.Range("L16").Value = 4666.3335
$ws.Range("M16").Value = -783
$ws.Range("N16").Value = -5006.3335
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H55").Value = 538.3158
$ws.Range("I55").Value = 523.6667
$ws.Range("K55").Value = 523.6667
$ws.Range("M55").Value = -350.6667
$ws.Range("H68").Value = 11602.556
$ws.Range("I68").Value = 980
$ws.Range("J68").Value = 12930.375
$ws.Range("K68").Value = 980
$ws.Range("L68").Value = 12930.375
$ws.Range("M68").Value = -231
$ws.Range("N68").Value = -14428.375
$ws.Range("H71").Value = 11602.556
$ws.Range("I71").Value = 980
$ws.Range("J71").Value = 12930.375
$ws.Range("K71").Value = 4900
$ws.Range("L71").Value = 64651.875
$ws.Range("M71").Value = -1156
$ws.Range("N71").Value = -72139.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 37552.5
$ws.Range("I70").Value = 35000
$ws.Range("J70").Value = 40105
$ws.Range("K70").Value = 35000
$ws.Range("L70").Value = 40105
$ws.Range("M70").Value = -34685
$ws.Range("N70").Value = -40735
$ws.Range("H73").Value = 37552.5
$ws.Range("I73").Value = 35000
$ws.Range("J73").Value = 40105
$ws.Range("K73").Value = 35000
$ws.Range("L73").Value = 40105
$ws.Range("M73").Value = -33908
$ws.Range("N73").Value = -42289
$ws.Range("H81").Value = 18190360
$ws.Range("I81").Value = 6322.6665
$ws.Range("J81").Value = 25009376
$ws.Range("K81").Value = 12645.333
$ws.Range("L81").Value = 50018752
$ws.Range("M81").Value = -11584.333
$ws.Range("N81").Value = -50020874
$ws.Range("H84").Value = 18190360
$ws.Range("I84").Value = 6322.6665
$ws.Range("J84").Value = 25009376
$ws.Range("K84").Value = 63226.665
$ws.Range("L84").Value = 250093760
$ws.Range("M84").Value = -57922.665
$ws.Range("N84").Value = -250104368
$ws.Range("H107").Value = 409.0625
$ws.Range("I107").Value = 350.07144
$ws.Range("K107").Value = 1050.21432
$ws.Range("M107").Value = 869.78568
$ws.Range("H122").Value = 1958.7812
$ws.Range("I122").Value = 2005.44
$ws.Range("J122").Value = 1792.1428
$ws.Range("K122").Value = 6016.32
$ws.Range("L122").Value = 5376.428400000001
$ws.Range("M122").Value = -3566.32
$ws.Range("N122").Value = -10276.4284
$ws.Range("H126").Value = 1780.1904
$ws.Range("I126").Value = 1634.4117
$ws.Range("K126").Value = 4903.2351
$ws.Range("M126").Value = -2433.2351
$ws.Range("H129").Value = 76666.664
$ws.Range("H132").Value = 468856.75
$ws.Range("I132").Value = 743590.3
$ws.Range("J132").Value = 5243.875
$ws.Range("K132").Value = 2230770.9
$ws.Range("L132").Value = 15731.625
$ws.Range("M132").Value = -2228240.9
$ws.Range("N132").Value = -20791.625

